# ============================================================================
# bdd/test.xlsx - add a new "crit_edad" criteria column (D) between
# preg_secc0 (C) and preg_test_1 (old D).  Existing D:G columns (and their
# header comments) shift right to E:H.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1. Re-anchor the header comments that sit to the right of the insertion
#    point (D:G) one column over (E:H), working from the far end inward so
#    nothing gets clobbered. C1's comment is untouched (it doesn't move).
# ----------------------------------------------------------------------
$t = $ws.Range("G1").Comment.Text()
$ws.Range("G1").Comment.Delete()
$ws.Range("H1").AddComment($t)

$t = $ws.Range("F1").Comment.Text()
$ws.Range("F1").Comment.Delete()
$ws.Range("G1").AddComment($t)

$t = $ws.Range("E1").Comment.Text()
$ws.Range("E1").Comment.Delete()
$ws.Range("F1").AddComment($t)

$t = $ws.Range("D1").Comment.Text()
$ws.Range("D1").Comment.Delete()
$ws.Range("E1").AddComment($t)

# ----------------------------------------------------------------------
# 2. Insert a new blank column at D - shifts values in D:G to E:H.
# ----------------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# ----------------------------------------------------------------------
# 3. Drop the now-orphaned "empty but styled" placeholder cells that used
#    to live in E:F and have shifted to F:G - they carried no real data.
# ----------------------------------------------------------------------
$ws.Range("F2:G7").ClearContents()

# ----------------------------------------------------------------------
# 4. Populate the new crit_edad column.
# ----------------------------------------------------------------------
$ws.Range("D1").Value = "crit_edad"
$ws.Range("D1").AddComment("Criterio de preg_edad, donde es 1 si la respuesta es <= 40")

$d2 = $ws.Range("D2")
$d2.Formula = '="1"'
$d2.Copy()
$d2.PasteSpecial(-4163)

$d3 = $ws.Range("D3")
$d3.Formula = '="0"'
$d3.Copy()
$d3.PasteSpecial(-4163)

